# Updates cryptocurrency price/volume data to reflect the latest scrape.
# Generated from the authoritative cell-by-cell diff of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.923.62'
$ws.Range("E2").Value = '  -1.03%  '
$ws.Range("D3").Value = '2.569.27'
$ws.Range("E3").Value = '  +2.25%  '
$ws.Range("D5").Value = '''302.60'
$ws.Range("E5").Value = '  +2.18%  '
$ws.Range("D6").Value = '''96.98'
$ws.Range("E6").Value = '  +2.94%  '
$ws.Range("E7").Value = '  +0.65%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("E11").Value = '  +0.97%  '
$ws.Range("D12").Value = '''7.64'
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("E13").Value = '  +6.72%  '
$ws.Range("D14").Value = '2.535.55'
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("E15").Value = '  +2.14%  '
$ws.Range("D16").Value = '''14.35'
$ws.Range("E16").Value = '  +2.25%  '
$ws.Range("D17").Value = '42.946.35'
$ws.Range("E17").Value = '  -0.94%  '
$ws.Range("D18").Value = '0.0₃0996'
$ws.Range("E18").Value = '  +3.43%  '
$ws.Range("D19").Value = '''12.89'
$ws.Range("E19").Value = '  +5.24%  '
$ws.Range("D20").Value = '''6.63'
$ws.Range("E20").Value = '  +1.53%  '
$ws.Range("D21").Value = '''71.99'
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").Value = '''254.21'
$ws.Range("E22").Value = '  -1.90%  '
$ws.Range("E23").Value = '  +2.25%  '
$ws.Range("E24").Value = '  -0.94%  '
$ws.Range("D25").Value = '''28.82'
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").Value = '''10.28'
$ws.Range("E27").Value = '  +2.80%  '
$ws.Range("D28").Value = '''37.45'
$ws.Range("E28").Value = '  +1.55%  '
$ws.Range("E29").Value = '  -5.21%  '
$ws.Range("D30").Value = '''6.03'
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("D31").Value = '''155.34'
$ws.Range("E31").Value = '  +3.56%  '
$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D32").Value = '''2.16'
$ws.Range("E32").Value = '  +1.21%  '
$ws.Range("B33").Value = 'WEMIXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D33").Value = '''2.76'
$ws.Range("E33").Value = '  -0.46%  '
$ws.Range("E34").Value = '  -1.79%  '
$ws.Range("D35").Value = '''0.0807'
$ws.Range("E35").Value = '  +1.45%  '
$ws.Range("D36").Value = '''18.34'
$ws.Range("E36").Value = '  +14.06%  '
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("E38").Value = '  +0.92%  '
$ws.Range("D39").Value = '''23.41'
$ws.Range("E39").Value = '  -0.59%  '
$ws.Range("E40").Value = '  -1.10%  '
$ws.Range("D41").Value = '''0.0310'
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '''3.87'
$ws.Range("E42").Value = '  +2.01%  '
$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").Value = '''2.05'
$ws.Range("E43").Value = '  +25.95%  '
$ws.Range("D44").Value = '2.070.88'
$ws.Range("E44").Value = '  +3.12%  '
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("D46").Value = '''9.23'
$ws.Range("E46").Value = '  +3.70%  '
$ws.Range("D47").Value = '''85.40'
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("D48").Value = '''76.60'
$ws.Range("E48").Value = '  +13.20%  '
$ws.Range("D49").Value = '''106.49'
$ws.Range("E49").Value = '  +3.52%  '
$ws.Range("D50").Value = '2.820.50'
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("E51").Value = '  +2.18%  '
